$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First, copy the A6 cell format (bold/border/centered) down to the two new rows ---
# (xlPasteFormats = -4122) so new rows keep the same "index" column style as the existing data rows.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)

# --- Update existing row 6 (job moved from IDEXX internship -> Deloitte Data Science Manager) ---
$ws.Range("C6").Value = 'Data Science Manager'
$ws.Range("D6").Value = 'Deloitte'
$ws.Range("E6").Value = 'do you have a passion for artificial intelligence, machine learning, and data analysis? do you yearn to have the impact of your work recognized and valued by more than just your development team? do you constantly wonder what you could build if only you had access to world-class data sets and computing resources?if yes, we have just the role for you.in deloitte''s audit and assurance business, we make businesses and markets better. an audit is more than an obligation; it is an opportunity to see further and deeper into businesses. in our role as independent auditors, we enhance trust in the companies we audit, helping a multitrillion dollar capital markets system function with greater confidence. as we aspire to the very highest standards of audit quality, we deliver deeper insights that can help clients become more effective organizations.deloitte''s audit and assurance business embraces the promise of artificial intelligence and machine learning, with various forms of ai and ml embedded in the audit technology solutions currently used by our 10,000+ practitioners today. within our dedicated data science organization, a subgroup of our award-winning audit transformation group, we continue that journey enabling the next generation of ai-enabled solutions that usher in a future that completely transforms the way our practitioners perform their audit work, and the insights we can provide to our clients.you will be joining a growing team of talented professionals in a fast-paced yet collaborative startup-like environment dedicated to realizing the deloitte audit and assurance vision of an ai-enabled audit. you will be leveraging the most advanced technologies in machine learning, natural language processing, time-series modeling, and reinforcement learning to lead our business into the future.as a data science manager, you will lead the technical and technological components of our data science workstreams and ai and ml solutions. you will work hand-in-hand with subject matter experts to ensure the inputs and outputs suit the intended user experience and audit workflow. you will manage data scientists, junior data scientists, and developers to complete the data science objectives that are required to fully develop and deploy a production-level ml solution. this includes performing or guiding research in an independent fashion, and interacting closely with key stakeholders with varying levels of machine learning experience.specifically, you will be expected to:have a profound understanding of the state of the art of a multitude of fields in artificial intelligence, including but not limited to nlp, probabilistic graphical models, time-series analysis, and weak supervised learning among otherslead the application of rigorous data science within your workstream, managing and supervising junior resources to do so; particularly:perform or lead the performance of exploratory data analysis to understand relationships and opportunities to influence outcomes, while being able to quickly iterate common feature transformation and model types to find the best predictive modelsdevelop and document proofs of concept to verify your ideas, including counterfactual explanations for interpretabilityclose the loop to make sure that the proposed solution is performing as it should and is correctly understoodcollaborate with subject matter experts to obtain a deep understanding of the underlying business problem, and to define and refine the corresponding technical solutionco-lead the planning and direction of a project with subject matter experts, and effectively prioritize goals and objectivesidentify opportunities to apply the latest advancements in machine learning and artificial intelligence to build, test, and validate predictive modelsmake impactful contributions to internal discussions on emerging machine learning methodologiesinfluence machine learning strategy for current and prospective workstreamsactively mentor data scientists and junior data scientists on good software practicesdevelop and embed automated processes for predictive model validation, deployment, and implementationarchitect ml pipelines and actively contribute high-quality, production-ready code (readable, well-tested, with well-designed apis)qualificationsrequired:phd in a quantitative field (computer science, engineering, mathematics, physics, machine learning, statistics)6+ years of industry experience leading the design, development, and deployment of machine learning modelsexperience being the technical lead for multiple project teams simultaneouslyprevious experience mentoring, training, and developing junior members of the team; experience in employee performance reviewsexpert understanding of python and other common languagesdeep understanding of machine learning model development life cyclesextensive experience using common machine learning and deep learning libraries and techniques, including tensorflow, pytorch, and big data platformsextensive experience with cloud-based ecosystems (azure, gcp, aws)experience in an agile working environment and related project management tools (jira, azure devops, etc.)demonstrated ability to write high-quality, production-ready code (readable, well-tested, well-documented, with well-designed apis)demonstrated ability to develop novel machine learning methods that go beyond putting together existing open-source code, and to apply problem-solving skills to complex issuesexperience with version control practices and tools (git, etc.)fluency in both structured and unstructured data (sql, nosql)solid understanding of docker, jenkins, kubernetes, and other devops toolsexcellent written and verbal communication skillsability to travel when necessarypreferred:10+ years of industry experience leading the design, development, and deployment of machine learning modelsprior scientific publication history. outstanding academic track record as evidenced by top tier publications.strong competency for additional coding languages (r, etc.)strong project management and delivery experience, including budget oversight and staffing of project teams including time managementextensive experience with microsoft azure, including certification in machine learningexperience with machine learning pipelines (azure ml)experience with ml ops and related governance processes, particularly within a regulated industrystrong presentation skills using microsoft office suite (visio, powerpoint, etc.)understanding of the capital markets, and the role public accounting firms'
$ws.Range("F6").Value = 'do you have a passion for artificial intelligence, machine learning, and data analysis? do you yearn to have the impact of your work recognized and valued by more than just your development team? do you constantly wonder what you could build if only you had access to world-class data sets and computing resources?if yes, we have just the role for you.in deloitte''s audit and assurance business, we make businesses and markets better. an audit is more than an obligation; it is an opportunity to see further and deeper into businesses. in our role as independent auditors, we enhance trust in the companies we audit, helping a multitrillion dollar capital markets system function with greater confidence. as we aspire to the very highest standards of audit quality, we deliver deeper insights that can help clients become more effective organizations.deloitte''s audit and assurance business embraces the promise of artificial intelligence and machine learning, with various forms of ai and ml embedded in the audit technology solutions currently used by our 10,000+ practitioners today. within our dedicated data science organization, a subgroup of our award-winning audit transformation group, we continue that journey enabling the next generation of ai-enabled solutions that usher in a future that completely transforms the way our practitioners perform their audit work, and the insights we can provide to our clients.you will be joining a growing team of talented professionals in a fast-paced yet collaborative startup-like environment dedicated to realizing the deloitte audit and assurance vision of an ai-enabled audit. you will be leveraging the most advanced technologies in machine learning, natural language processing, time-series modeling, and reinforcement learning to lead our business into the future.as a data science manager, you will lead the technical and technological components of our data science workstreams and ai and ml solutions. you will work hand-in-hand with subject matter experts to ensure the inputs and outputs suit the intended user experience and audit workflow. you will manage data scientists, junior data scientists, and developers to complete the data science objectives that are required to fully develop and deploy a production-level ml solution. this includes performing or guiding research in an independent fashion, and interacting closely with key stakeholders with varying levels of machine learning experience.specifically, you will be expected to:have a profound understanding of the state of the art of a multitude of fields in artificial intelligence, including but not limited to nlp, probabilistic graphical models, time-series analysis, and weak supervised learning among otherslead the application of rigorous data science within your workstream, managing and supervising junior resources to do so; particularly:perform or lead the performance of exploratory data analysis to understand relationships and opportunities to influence outcomes, while being able to quickly iterate common feature transformation and model types to find the best predictive modelsdevelop and document proofs of concept to verify your ideas, including counterfactual explanations for interpretabilityclose the loop to make sure that the proposed solution is performing as it should and is correctly understoodcollaborate with subject matter experts to obtain a deep understanding of the underlying business problem, and to define and refine the corresponding technical solutionco-lead the planning and direction of a project with subject matter experts, and effectively prioritize goals and objectivesidentify opportunities to apply the latest advancements in machine learning and artificial intelligence to build, test, and validate predictive modelsmake impactful contributions to internal discussions on emerging machine learning methodologiesinfluence machine learning strategy for current and prospective workstreamsactively mentor data scientists and junior data scientists on good software practicesdevelop and embed automated processes for predictive model validation, deployment, and implementationarchitect ml pipelines and actively contribute high-quality, production-ready code (readable, well-tested, with well-designed apis)qualificationsrequired:phd in a quantitative field (computer science, engineering, mathematics, physics, machine learning, statistics)6+ years of industry experience leading the design, development, and deployment of machine learning modelsexperience being the technical lead for multiple project teams simultaneouslyprevious experience mentoring, training, and developing junior members of the team; experience in employee performance reviewsexpert understanding of python and other common languagesdeep understanding of machine learning model development life cyclesextensive experience using common machine learning and deep learning libraries and techniques, including tensorflow, pytorch, and big data platformsextensive experience with cloud-based ecosystems (azure, gcp, aws)experience in an agile working environment and related project management tools (jira, azure devops, etc.)demonstrated ability to write high-quality, production-ready code (readable, well-tested, well-documented, with well-designed apis)demonstrated ability to develop novel machine learning methods that go beyond putting together existing open-source code, and to apply problem-solving skills to complex issuesexperience with version control practices and tools (git, etc.)fluency in both structured and unstructured data (sql, nosql)solid understanding of docker, jenkins, kubernetes, and other devops toolsexcellent written and verbal communication skillsability to travel when necessarypreferred:10+ years of industry experience leading the design, development, and deployment of machine learning modelsprior scientific publication history. outstanding academic track record as evidenced by top tier publications.strong competency for additional coding languages (r, etc.)strong project management and delivery experience, including budget oversight and staffing of project teams including time managementextensive experience with microsoft azure, including certification in machine learningexperience with machine learning pipelines (azure ml)experience with ml ops and related governance processes, particularly within a regulated industrystrong presentation skills using microsoft office suite (visio, powerpoint, etc.)understanding of the capital markets, and the role public accounting firms'
$ws.Range("G6").Value = 0.09478021978021978
$ws.Range("H6").Value = 0.4748473748473748
$ws.Range("I6").Value = '[''include'', ''have'', ''learn'', ''lead'', ''perform'', ''use'', ''enable'', ''develop'', ''build'', ''make'', ''help'', ''work'', ''manage'', ''apply'', ''test'', ''design'', ''write'', ''yearn'', ''recognize'', ''value'', ''wonder'', ''see'', ''enhance'', ''audit'', ''aspire'', ''deliver'', ''become'', ''embrace'', ''embed'', ''win'', ''continue'', ''transform'', ''provide'', ''join'', ''grow'', ''pace'', ''dedicate'', ''realize'', ''leverage'', ''workstream'', ''ai'', ''ensure'', ''suit'', ''intend'', ''complete'', ''require'', ''deploy'', ''guide'', ''interact'', ''vary'', ''expect'', ''limit'', ''supervise'', ''do'', ''understand'', ''influence'', ''iterate'', ''find'', ''verify'', ''interpretabilityclose'', ''propose'', ''obtain'', ''underlie'', ''define'', ''refine'', ''correspond'', ''prioritize'', ''emerge'', ''embe'', ''automate'', ''implementationarchitect'', ''contribute'', ''mentor'', ''base'', ''etc.)demonstrate'', ''document'', ''apis)demonstrate'', ''go'', ''put'', ''exist'', ''solve'', ''travel'', ''evidence'', ''code'', ''azure'']'

# --- Add new row 7 (Entry Level Data Scientist / SynergisticIT) ---
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 'https://www.linkedin.com/jobs/view/3036847597'
$ws.Range("C7").Value = 'Entry Level Data Scientist'
$ws.Range("D7").Value = 'SynergisticIT'
$ws.Range("E7").Value = 'at synergisticit, we aim to bring aboard it professionals to help them build a rewarding career in cutting-edge technologies. being in the industry for more than 10 years, we provide a splendid range of lucrative opportunities to sustain a position in our top tech clients like google, apple, cognizant, client, paypal, to name a few. our seasoned team firmly believes that the new tech talent can scale any business if given the right opportunity. we value your integrity, hard work, and commitment to make a difference in the technical sphere. for this reason, we focus on providing end-to-end career assistance and enhancing your already existing it skills and knowledge. currently, we are looking for qualified entry-level data scientists who can apply data science principles to design, test, implement, and develop data-based solutions, including reporting, auditing, and preparing large databases for statistical analysis. minimum background and qualifications requirement bachelor''s degree or master''s degree in computer engineering, computer science, mathematics, electrical engineering, information systems, or it must have mathematics or statistics background technical and soft skills required experience in python programming and understanding of the software development life cycle. knowledge of linear algebra, statistics, and mathematics concepts. excellent written and verbal communication skills. highly motivated, self-learner, team player, and technically inquisitive. strong work ethics and creative problem-solving abilities. preferred skills deep learning data visualization nlp scala django roles and responsibilities collaborate with dynamic teams of engineers, developers, and scientists who research and integrate algorithms to develop an application, software, and computer system solutions to address complex data problems. assess project requirements and develop data analysis algorithms. engage developers to share their opinions, knowledge, and recommendations to meet the deliverables. contribute to technical solutions and implement software analyses to unlock the secrets held by big data sets. integrate components like web-based ui, commercial indexing products, and access control mechanisms to create operational information and knowledge discovery systems. benefits competitive salary flexible work schedule part-time off e-verified no relocation h1b filing on job technical support skill enhancement opportunity to work with fortune 500 companies who should apply? recent it graduates looking to build a solid career in the tech industry. if you''re lured by the endless possibilities presented by ai, machine learning, iot, and data science, this job opportunity can be the right career path for you. candidate''s outcome httpswww.synergisticit.comcandidate-outcomes httpswww.synergisticit.comcandidate-outcomes no third-party candidates or c2c candidates if you are interested, please apply to the posting. no phone calls please, shortlisted candidates would be reached out.'
$ws.Range("F7").Value = 'at synergisticit, we aim to bring aboard it professionals to help them build a rewarding career in cutting-edge technologies. being in the industry for more than 10 years, we provide a splendid range of lucrative opportunities to sustain a position in our top tech clients like google, apple, cognizant, client, paypal, to name a few. our seasoned team firmly believes that the new tech talent can scale any business if given the right opportunity. we value your integrity, hard work, and commitment to make a difference in the technical sphere. for this reason, we focus on providing end-to-end career assistance and enhancing your already existing it skills and knowledge. currently, we are looking for qualified entry-level data scientists who can apply data science principles to design, test, implement, and develop data-based solutions, including reporting, auditing, and preparing large databases for statistical analysis. minimum background and qualifications requirement bachelor''s degree or master''s degree in computer engineering, computer science, mathematics, electrical engineering, information systems, or it must have mathematics or statistics background technical and soft skills required experience in python programming and understanding of the software development life cycle. knowledge of linear algebra, statistics, and mathematics concepts. excellent written and verbal communication skills. highly motivated, self-learner, team player, and technically inquisitive. strong work ethics and creative problem-solving abilities. preferred skills deep learning data visualization nlp scala django roles and responsibilities collaborate with dynamic teams of engineers, developers, and scientists who research and integrate algorithms to develop an application, software, and computer system solutions to address complex data problems. assess project requirements and develop data analysis algorithms. engage developers to share their opinions, knowledge, and recommendations to meet the deliverables. contribute to technical solutions and implement software analyses to unlock the secrets held by big data sets. integrate components like web-based ui, commercial indexing products, and access control mechanisms to create operational information and knowledge discovery systems. benefits competitive salary flexible work schedule part-time off e-verified no relocation h1b filing on job technical support skill enhancement opportunity to work with fortune 500 companies who should apply? recent it graduates looking to build a solid career in the tech industry. if you''re lured by the endless possibilities presented by ai, machine learning, iot, and data science, this job opportunity can be the right career path for you. candidate''s outcome httpswww.synergisticit.comcandidate-outcomes httpswww.synergisticit.comcandidate-outcomes no third-party candidates or c2c candidates if you are interested, please apply to the posting. no phone calls please, shortlisted candidates would be reached out.'
$ws.Range("G7").Value = 0.1510024350649351
$ws.Range("H7").Value = 0.4198816287878787
$ws.Range("I7").Value = '[''apply'', ''develop'', ''build'', ''provide'', ''look'', ''implement'', ''base'', ''integrate'', ''aim'', ''bring'', ''help'', ''cut'', ''sustain'', ''name'', ''believe'', ''scale'', ''give'', ''value'', ''make'', ''focus'', ''enhance'', ''exist'', ''design'', ''include'', ''prepare'', ''have'', ''require'', ''write'', ''solve'', ''learn'', ''collaborate'', ''research'', ''address'', ''assess'', ''engage'', ''share'', ''meet'', ''contribute'', ''unlock'', ''hold'', ''create'', ''benefit'', ''-'', ''verify'', ''file'', ''work'', ''graduate'', ''lure'', ''present'', ''please'', ''shortlist'', ''reach'']'

# --- Add new row 8 (Entry Level Data Scientist / Dice) ---
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 'https://www.linkedin.com/jobs/view/3041899464'
$ws.Range("C8").Value = 'Entry Level Data Scientist'
$ws.Range("D8").Value = 'Dice'
$ws.Range("E8").Value = 'dice is the leading career destination for tech experts at every stage of their careers. our client, synergisticit, is seeking the following. apply via dice today!at synergisticit, we aim to bring aboard it professionals to help them build a rewarding career in cutting-edge technologies. being in the industry for more than 10 years, we provide a splendid range of lucrative opportunities to sustain a position in our top tech clients like google, apple, cognizant, client, paypal, to name a few. our seasoned team firmly believes that the new tech talent can scale any business if given the right opportunity. we value your integrity, hard work, and commitment to make a difference in the technical sphere. for this reason, we focus on providing end-to-end career assistance and enhancing your already existing it skills and knowledge. currently, we are looking for qualified entry-level data scientists who can apply data science principles to design, test, implement, and develop data-based solutions, including reporting, auditing, and preparing large databases for statistical analysis. minimum background and qualifications requirement bachelor''s degree or master''s degree in computer engineering, computer science, mathematics, electrical engineering, information systems, or it must have mathematics or statistics background technical and soft skills required experience in python programming and understanding of the software development life cycle. knowledge of linear algebra, statistics, and mathematics concepts. excellent written and verbal communication skills. highly motivated, self-learner, team player, and technically inquisitive. strong work ethics and creative problem-solving abilities. preferred skills deep learning data visualization nlp scala django roles and responsibilities collaborate with dynamic teams of engineers, developers, and scientists who research and integrate algorithms to develop an application, software, and computer system solutions to address complex data problems. assess project requirements and develop data analysis algorithms. engage developers to share their opinions, knowledge, and recommendations to meet the deliverables. contribute to technical solutions and implement software analyses to unlock the secrets held by big data sets. integrate components like web-based ui, commercial indexing products, and access control mechanisms to create operational information and knowledge discovery systems. benefits competitive salary flexible work schedule part-time off e-verified no relocation h1b filing on job technical support skill enhancement opportunity to work with fortune 500 companies who should apply? recent it graduates looking to build a solid career in the tech industry. if you''re lured by the endless possibilities presented by ai, machine learning, iot, and data science, this job opportunity can be the right career path for you. candidate''s outcome httpswww.synergisticit.comcandidate-outcomes httpswww.synergisticit.comcandidate-outcomes no third-party candidates or c2c candidates if you are interested, please apply to the posting. no phone calls please, shortlisted candidates would be reached out.'
$ws.Range("F8").Value = 'dice is the leading career destination for tech experts at every stage of their careers. our client, synergisticit, is seeking the following. apply via dice today!at synergisticit, we aim to bring aboard it professionals to help them build a rewarding career in cutting-edge technologies. being in the industry for more than 10 years, we provide a splendid range of lucrative opportunities to sustain a position in our top tech clients like google, apple, cognizant, client, paypal, to name a few. our seasoned team firmly believes that the new tech talent can scale any business if given the right opportunity. we value your integrity, hard work, and commitment to make a difference in the technical sphere. for this reason, we focus on providing end-to-end career assistance and enhancing your already existing it skills and knowledge. currently, we are looking for qualified entry-level data scientists who can apply data science principles to design, test, implement, and develop data-based solutions, including reporting, auditing, and preparing large databases for statistical analysis. minimum background and qualifications requirement bachelor''s degree or master''s degree in computer engineering, computer science, mathematics, electrical engineering, information systems, or it must have mathematics or statistics background technical and soft skills required experience in python programming and understanding of the software development life cycle. knowledge of linear algebra, statistics, and mathematics concepts. excellent written and verbal communication skills. highly motivated, self-learner, team player, and technically inquisitive. strong work ethics and creative problem-solving abilities. preferred skills deep learning data visualization nlp scala django roles and responsibilities collaborate with dynamic teams of engineers, developers, and scientists who research and integrate algorithms to develop an application, software, and computer system solutions to address complex data problems. assess project requirements and develop data analysis algorithms. engage developers to share their opinions, knowledge, and recommendations to meet the deliverables. contribute to technical solutions and implement software analyses to unlock the secrets held by big data sets. integrate components like web-based ui, commercial indexing products, and access control mechanisms to create operational information and knowledge discovery systems. benefits competitive salary flexible work schedule part-time off e-verified no relocation h1b filing on job technical support skill enhancement opportunity to work with fortune 500 companies who should apply? recent it graduates looking to build a solid career in the tech industry. if you''re lured by the endless possibilities presented by ai, machine learning, iot, and data science, this job opportunity can be the right career path for you. candidate''s outcome httpswww.synergisticit.comcandidate-outcomes httpswww.synergisticit.comcandidate-outcomes no third-party candidates or c2c candidates if you are interested, please apply to the posting. no phone calls please, shortlisted candidates would be reached out.'
$ws.Range("G8").Value = 0.146426603699331
$ws.Range("H8").Value = 0.410188246097337
$ws.Range("I8").Value = '[''apply'', ''develop'', ''build'', ''provide'', ''look'', ''implement'', ''base'', ''integrate'', ''lead'', ''seek'', ''aim'', ''bring'', ''help'', ''cut'', ''sustain'', ''name'', ''believe'', ''scale'', ''give'', ''value'', ''make'', ''focus'', ''enhance'', ''exist'', ''design'', ''include'', ''prepare'', ''have'', ''require'', ''write'', ''solve'', ''learn'', ''collaborate'', ''research'', ''address'', ''assess'', ''engage'', ''share'', ''meet'', ''contribute'', ''unlock'', ''hold'', ''create'', ''benefit'', ''-'', ''verify'', ''file'', ''work'', ''graduate'', ''lure'', ''present'', ''please'', ''shortlist'', ''reach'']'

